# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos" edit:
# rotates the three named-worker rows (16-18), reverses the "MIGUEL GONZALEZ" period
# list (rows 19-33), and updates the corresponding Valor Mora / Salario Basico figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 16-18: rotate worker identity (doc #, name) one position up, wrapping ---
$ws.Range("C16").Value = "1238338053"
$ws.Range("D16").Value = "DANIEL EDUARDO BERRIO BELTRAN"

$ws.Range("C17").Value = "1050969229"
$ws.Range("D17").Value = "VICTOR VICTOR RAMOS JIMENEZ"

$ws.Range("C18").Value = "1047373924"
$ws.Range("D18").Value = "ALEJANDRO FRANCO BARRIOS"

# --- Rows 16-18: updated Valor Mora (F) / Salario Basico (G) ---
$ws.Range("F16").Value = 121333
$ws.Range("G16").Value = 1100000

$ws.Range("F17").Value = 89600
$ws.Range("G17").Value = 11200000

$ws.Range("F18").Value = 177333
$ws.Range("G18").Value = 13300000

# --- Rows 19-33: reverse the "Periodo Mora" sequence for MIGUEL GONZALEZ ---
$ws.Range("E19").Value = "2309"
$ws.Range("E20").Value = "2308"
$ws.Range("E21").Value = "2307"
$ws.Range("E22").Value = "2306"
$ws.Range("E23").Value = "2305"
$ws.Range("E24").Value = "2304"
$ws.Range("E25").Value = "2303"
$ws.Range("E26").Value = "2302"
$ws.Range("E27").Value = "2301"
$ws.Range("E28").Value = "2212"
$ws.Range("E29").Value = "2211"
$ws.Range("E30").Value = "2210"
$ws.Range("E31").Value = "2209"
$ws.Range("E32").Value = "2208"
$ws.Range("E33").Value = "2207"

# --- Valor Mora swap that tracks the Periodo Mora reversal on rows 19 & 33 ---
$ws.Range("F19").Value = 252720
$ws.Range("F33").Value = 98280
